# Apply "new progress as of date 04 Nov 2025" update to the Training Dashboard sheet:
# - Column H (PERIOD TO EXPIRE) decreases by 1 for each data row (3-21)
# - Column I (LAST UPDATE) changes from "03-Nov-2025" to "04-Nov-2025" for each data row (3-21)

$ws = $excel.ActiveWorkbook.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 21; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value2 = "04-Nov-2025"
}
